$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# Data rows (force numeric-looking strings to remain text)
$ws.Range("A2").Value = " Dubai (DSC)"
$ws.Range("B2").Value = " October 27 2020"
$ws.Range("C2").Value = "Sunrisers won by 88 runs"
$ws.Range("D2").Value = "Delhi Capitals"
$ws.Range("E2").Value = "Sunrisers Hyderabad"
$ws.Range("F2").Value = "Anrich Nortje "
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "1"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "3"
$ws.Range("H2").ClearFormats()
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "0"
$ws.Range("I2").ClearFormats()
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "0"
$ws.Range("J2").ClearFormats()
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "33.33"
$ws.Range("K2").ClearFormats()
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " November 05 2020"
$ws.Range("C3").Value = "Mumbai won by 57 runs"
$ws.Range("D3").Value = "Delhi Capitals"
$ws.Range("E3").Value = "Mumbai Indians"
$ws.Range("F3").Value = "Anrich Nortje "
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "0"
$ws.Range("H3").ClearFormats()
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "0"
$ws.Range("I3").ClearFormats()
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "0"
$ws.Range("J3").ClearFormats()
$ws.Range("K3").Value = "-"
$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " September 20 2020"
$ws.Range("C4").Value = "Match tied (Capitals won the one-over eliminator)"
$ws.Range("D4").Value = "Delhi Capitals"
$ws.Range("E4").Value = "Kings XI Punjab"
$ws.Range("F4").Value = "Anrich Nortje "
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "3"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "1"
$ws.Range("H4").ClearFormats()
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "0"
$ws.Range("I4").ClearFormats()
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "0"
$ws.Range("J4").ClearFormats()
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "300.00"
$ws.Range("K4").ClearFormats()
$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " October 24 2020"
$ws.Range("C5").Value = "KKR won by 59 runs"
$ws.Range("D5").Value = "Delhi Capitals"
$ws.Range("E5").Value = "Kolkata Knight Riders"
$ws.Range("F5").Value = "Anrich Nortje "
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0"
$ws.Range("G5").ClearFormats()
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "0"
$ws.Range("H5").ClearFormats()
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "0"
$ws.Range("I5").ClearFormats()
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "0"
$ws.Range("J5").ClearFormats()
$ws.Range("K5").Value = "-"
$ws.Range("A6").Value = " Abu Dhabi"
$ws.Range("B6").Value = " September 29 2020"
$ws.Range("C6").Value = "Sunrisers won by 15 runs"
$ws.Range("D6").Value = "Delhi Capitals"
$ws.Range("E6").Value = "Sunrisers Hyderabad"
$ws.Range("F6").Value = "Anrich Nortje "
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "3"
$ws.Range("G6").ClearFormats()
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "2"
$ws.Range("H6").ClearFormats()
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "0"
$ws.Range("I6").ClearFormats()
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "0"
$ws.Range("J6").ClearFormats()
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "150.00"
$ws.Range("K6").ClearFormats()
